$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8.350541
$ws.Range("H2").Value = 25.051623
$ws.Range("I2").Value = 0.3086892463293835
$ws.Range("J2").Value = 0.3086892463293835
$ws.Range("M2").Value = 11.188041
$ws.Range("N2").Value = 33.564123
$ws.Range("O2").Value = 0.1395103797998223
$ws.Range("P2").Value = 0.1395103797998223
$ws.Range("Q2").Value = 93.426195080181
$ws.Range("R2").Value = 840.8357557216291
$ws.Range("S2").Value = 0.0430653539955332
$ws.Range("T2").Value = 0.04306535399553321

$ws.Range("G3").Value = 8.350541
$ws.Range("H3").Value = 25.051623
$ws.Range("I3").Value = 0.3086892463293835
$ws.Range("J3").Value = 0.3086892463293835
$ws.Range("O3").Value = 0.4168441980730721
$ws.Range("P3").Value = 0.4168441980730722
$ws.Range("Q3").Value = 279.1488878683854
$ws.Range("R3").Value = 2512.339990815468
$ws.Range("S3").Value = 0.1286753213399529
$ws.Range("T3").Value = 0.1286753213399529

$ws.Range("G4").Value = 8.350541
$ws.Range("H4").Value = 25.051623
$ws.Range("I4").Value = 0.3086892463293835
$ws.Range("J4").Value = 0.3086892463293835
$ws.Range("M4").Value = 31.78201566666667
$ws.Range("N4").Value = 95.346047
$ws.Range("O4").Value = 0.3963089763847458
$ws.Range("P4").Value = 0.3963089763847459
$ws.Range("Q4").Value = 265.3970248871423
$ws.Range("R4").Value = 2388.573223984281
$ws.Range("S4").Value = 0.1223363192337766
$ws.Range("T4").Value = 0.1223363192337767

$ws.Range("G5").Value = 8.350541
$ws.Range("H5").Value = 25.051623
$ws.Range("I5").Value = 0.3086892463293835
$ws.Range("J5").Value = 0.3086892463293835
$ws.Range("M5").Value = 3.796148333333333
$ws.Range("N5").Value = 11.388445
$ws.Range("O5").Value = 0.04733644574235969
$ws.Range("P5").Value = 0.04733644574235969
$ws.Range("Q5").Value = 31.69989229958166
$ws.Range("R5").Value = 285.299030696235
$ws.Range("S5").Value = 0.01461225176012077
$ws.Range("T5").Value = 0.01461225176012077

$ws.Range("I6").Value = 0.4146406124520329
$ws.Range("J6").Value = 0.4146406124520329
$ws.Range("M6").Value = 11.188041
$ws.Range("N6").Value = 33.564123
$ws.Range("O6").Value = 0.1395103797998223
$ws.Range("P6").Value = 0.1395103797998223
$ws.Range("Q6").Value = 125.492854732536
$ws.Range("R6").Value = 1129.435692592824
$ws.Range("S6").Value = 0.05784666932361403
$ws.Range("T6").Value = 0.05784666932361403

$ws.Range("I7").Value = 0.4146406124520329
$ws.Range("J7").Value = 0.4146406124520329
$ws.Range("O7").Value = 0.4168441980730721
$ws.Range("P7").Value = 0.4168441980730722
$ws.Range("S7").Value = 0.1728405335860951
$ws.Range("T7").Value = 0.1728405335860951

$ws.Range("I8").Value = 0.4146406124520329
$ws.Range("J8").Value = 0.4146406124520329
$ws.Range("M8").Value = 31.78201566666667
$ws.Range("N8").Value = 95.346047
$ws.Range("O8").Value = 0.3963089763847458
$ws.Range("P8").Value = 0.3963089763847459
$ws.Range("Q8").Value = 356.4892080002373
$ws.Range("R8").Value = 3208.402872002136
$ws.Range("S8").Value = 0.1643257966884092
$ws.Range("T8").Value = 0.1643257966884093

$ws.Range("I9").Value = 0.4146406124520329
$ws.Range("J9").Value = 0.4146406124520329
$ws.Range("M9").Value = 3.796148333333333
$ws.Range("N9").Value = 11.388445
$ws.Range("O9").Value = 0.04733644574235969
$ws.Range("P9").Value = 0.04733644574235969
$ws.Range("Q9").Value = 42.58024182590666
$ws.Range("R9").Value = 383.2221764331599
$ws.Range("S9").Value = 0.01962761285391445
$ws.Range("T9").Value = 0.01962761285391445

$ws.Range("G10").Value = 7.484372666666666
$ws.Range("H10").Value = 22.453118
$ws.Range("I10").Value = 0.2766701412185836
$ws.Range("J10").Value = 0.2766701412185836
$ws.Range("M10").Value = 11.188041
$ws.Range("N10").Value = 33.564123
$ws.Range("O10").Value = 0.1395103797998223
$ws.Range("P10").Value = 0.1395103797998223
$ws.Range("Q10").Value = 83.735468253946
$ws.Range("R10").Value = 753.6192142855141
$ws.Range("S10").Value = 0.03859835648067506
$ws.Range("T10").Value = 0.03859835648067507

$ws.Range("G11").Value = 7.484372666666666
$ws.Range("H11").Value = 22.453118
$ws.Range("I11").Value = 0.2766701412185836
$ws.Range("J11").Value = 0.2766701412185836
$ws.Range("O11").Value = 0.4168441980730721
$ws.Range("P11").Value = 0.4168441980730722
$ws.Range("Q11").Value = 250.1938863952098
$ws.Range("R11").Value = 2251.744977556888
$ws.Range("S11").Value = 0.1153283431470241
$ws.Range("T11").Value = 0.1153283431470241

$ws.Range("G12").Value = 7.484372666666666
$ws.Range("H12").Value = 22.453118
$ws.Range("I12").Value = 0.2766701412185836
$ws.Range("J12").Value = 0.2766701412185836
$ws.Range("M12").Value = 31.78201566666667
$ws.Range("N12").Value = 95.346047
$ws.Range("O12").Value = 0.3963089763847458
$ws.Range("P12").Value = 0.3963089763847459
$ws.Range("Q12").Value = 237.8684493471718
$ws.Range("R12").Value = 2140.816044124546
$ws.Range("S12").Value = 0.1096468604625599
$ws.Range("T12").Value = 0.10964686046256

$ws.Range("G13").Value = 7.484372666666666
$ws.Range("H13").Value = 22.453118
$ws.Range("I13").Value = 0.2766701412185836
$ws.Range("J13").Value = 0.2766701412185836
$ws.Range("M13").Value = 3.796148333333333
$ws.Range("N13").Value = 11.388445
$ws.Range("O13").Value = 0.04733644574235969
$ws.Range("P13").Value = 0.04733644574235969
$ws.Range("Q13").Value = 28.41178882461222
$ws.Range("R13").Value = 255.70609942151
$ws.Range("S13").Value = 0.01309658112832447
$ws.Range("T13").Value = 0.01309658112832448
